# Update cost and robot comparison table
# - "Dollar to euro" sheet: add a new column I with the euro-converted
#   price for every part row (G*0.86), mirroring the existing total-only
#   conversion that used to live in G14.
# - "Euro suppliers" sheet: add a spacer row (row 13) between the two
#   cost tables, matching the layout of the "Dollar to euro" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Dollar to euro" - per-row euro price column (I2:I12)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Dollar to euro")

# Clarify that these headers refer to the dollar-denominated figures, now
# that a euro-converted column is being added alongside them.
$ws1.Range("E1").Value = "total cost (dollar):"
$ws1.Range("F1").Value = "Cost for 20 (dollar):"
$ws1.Range("G1").Value = "Total cost (dollar):"

for ($row = 2; $row -le 11; $row++) {
    $ws1.Range("I$row").Formula = "=G$row*0.86"
}
# Row 12 holds the column totals (G12 = SUM(G2:G11)); mirror that in I12.
$ws1.Range("I12").Formula = "=G12*0.86"

# Restore the selection/scroll position recorded after the edit.
$ws1.Range("J13").Select()

# ---------------------------------------------------------------------
# Sheet 2: "Euro suppliers" - blank spacer row between the two tables
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Euro suppliers")

$ws2.Range("A13:B13").NumberFormat = "General"
$ws2.Range("E13:H13").NumberFormat = "General"

$wb.Save()
